$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (row 1): extend the 0..13 sequence with two more columns ---
$ws.Range("P1").Value = 14
$ws.Range("P1").Borders.LineStyle = 1
$ws.Range("P1").Font.Bold = $true
$ws.Range("P1").HorizontalAlignment = -4108
$ws.Range("P1").VerticalAlignment = -4160

$ws.Range("Q1").Value = 15
$ws.Range("Q1").Borders.LineStyle = 1
$ws.Range("Q1").Font.Bold = $true
$ws.Range("Q1").HorizontalAlignment = -4108
$ws.Range("Q1").VerticalAlignment = -4160

# --- Data rows (2..25): swap I/K and M/O values, then append P and Q = 2 ---
for ($r = 2; $r -le 25; $r++) {
    $ws.Range("I$r").Value = 2
    $ws.Range("K$r").Value = 1
    $ws.Range("M$r").Value = 2
    $ws.Range("O$r").Value = 1
    $ws.Range("P$r").Value = 2
    $ws.Range("Q$r").Value = 2
}
